$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-15 Friday" "2023-12-16 Saturday"

Replace-Text "39×30=1170" "45×88=3960"
Replace-Text "15×87=1305" "62×92=5704"
Replace-Text "29×59=1711" "30×92=2760"
Replace-Text "48×70=3360" "45×21=945"
Replace-Text "90×74=6660" "67×77=5159"

Replace-Text "75×26=1950" "72×87=6264"
Replace-Text "17×82=1394" "26×36=936"
Replace-Text "28×22=616" "13×70=910"
Replace-Text "80×47=3760" "45×32=1440"
Replace-Text "35×16=560" "66×28=1848"

Replace-Text "99×33=3267" "63×83=5229"
Replace-Text "94×79=7426" "20×48=960"
Replace-Text "99×12=1188" "99×62=6138"
Replace-Text "52×97=5044" "84×87=7308"
Replace-Text "13×68=884" "46×96=4416"

Replace-Text "32×35=1120" "71×83=5893"
Replace-Text "14×62=868" "24×44=1056"
Replace-Text "91×90=8190" "35×95=3325"
Replace-Text "11×34=374" "62×91=5642"
Replace-Text "22×17=374" "16×13=208"

Replace-Text "87×72=6264" "80×76=6080"
Replace-Text "11×31=341" "24×34=816"
Replace-Text "69×95=6555" "26×20=520"
Replace-Text "57×56=3192" "14×72=1008"
Replace-Text "99×29=2871" "66×74=4884"
